$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: cell edits made before the filter/slicer selection changes -------
# "Fecha Actualización" (F) bumped for several rows. Writing a cell value in
# this sheet forces a row-height recompute, so each write is immediately
# followed by a restore of the row's original display height (rows that had
# no explicit height go through AutoFit to drop back to the default; rows
# that had an explicit wrap-text height get that height re-applied).
$ws.Range("F27").Value = 44147
$ws.Range("F27").EntireRow.AutoFit()

$ws.Range("F35").Value = 44147
$ws.Range("F35").EntireRow.RowHeight = 30

$ws.Range("F38").Value = 44147
$ws.Range("F38").EntireRow.AutoFit()

$ws.Range("F39").Value = 44147
$ws.Range("F39").EntireRow.AutoFit()

$ws.Range("F42").Value = 44147
$ws.Range("F42").EntireRow.RowHeight = 30

$ws.Range("F44").Value = 44147
$ws.Range("F44").EntireRow.RowHeight = 60

$ws.Range("F45").Value = 44147
$ws.Range("F45").EntireRow.AutoFit()

$ws.Range("F46").Value = 44147
$ws.Range("F46").EntireRow.RowHeight = 45

$ws.Range("F47").Value = 44147
$ws.Range("F47").EntireRow.RowHeight = 30

$ws.Range("F48").Value = 44147
$ws.Range("F48").EntireRow.RowHeight = 30

# Row 43's "Situación" (E) is reclassified from "Trabajando" to "Listo"
$ws.Range("E43").Value = "Listo"
$ws.Range("E43").EntireRow.RowHeight = 30

# --- Step 2: change the "Situación" slicer/filter selection ------------------
# from "Revisión" to "Trabajando" on the Tabla1 table (Situación is table column 4)
$lo = $ws.ListObjects.Item("Tabla1")
$lo.Range.AutoFilter(4, @("Trabajando"), 7) | Out-Null

# --- Step 3: leave the cursor where the editing user ended up ----------------
$ws.Range("F44").Select() | Out-Null
